$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(3, 6).Value = 234
$ws.Cells.Item(4, 6).Value = 195
$ws.Cells.Item(5, 6).Value = 1748
$ws.Cells.Item(6, 6).Value = 648
$ws.Cells.Item(8, 6).Value = 419
$ws.Cells.Item(9, 6).Value = 4133
$ws.Cells.Item(11, 6).Value = 438
$ws.Cells.Item(13, 6).Value = 973
$ws.Cells.Item(14, 6).Value = 1266
$ws.Cells.Item(17, 6).Value = 2929
$ws.Cells.Item(18, 6).Value = 1757
$ws.Cells.Item(20, 6).Value = 36
$ws.Cells.Item(21, 6).Value = 157
$ws.Cells.Item(23, 6).Value = 910
$ws.Cells.Item(24, 6).Value = 280
$ws.Cells.Item(25, 6).Value = 20
$ws.Cells.Item(26, 6).Value = 2195
$ws.Cells.Item(28, 6).Value = 2273
$ws.Cells.Item(29, 6).Value = 236
$ws.Cells.Item(30, 6).Value = 667
$ws.Cells.Item(31, 6).Value = 445
$ws.Cells.Item(33, 6).Value = 870
$ws.Cells.Item(34, 6).Value = 391
$ws.Cells.Item(35, 6).Value = 1043
$ws.Cells.Item(36, 6).Value = 864
$ws.Cells.Item(37, 6).Value = 1133
$ws.Cells.Item(38, 6).Value = 302
$ws.Cells.Item(39, 6).Value = 495
$ws.Cells.Item(40, 6).Value = 345
$ws.Cells.Item(42, 6).Value = 3441

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(10, 6).Value = 869

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(3, 6).Value = 234
$ws.Cells.Item(4, 6).Value = 195
$ws.Cells.Item(6, 6).Value = 1748
$ws.Cells.Item(7, 6).Value = 648
$ws.Cells.Item(9, 6).Value = 419
$ws.Cells.Item(10, 6).Value = 4133
$ws.Cells.Item(15, 6).Value = 1266
$ws.Cells.Item(16, 6).Value = 2929
$ws.Cells.Item(16, 7).Value = 68
$ws.Cells.Item(18, 6).Value = 1757
$ws.Cells.Item(20, 6).Value = 36
$ws.Cells.Item(22, 6).Value = 157
$ws.Cells.Item(23, 6).Value = 869
$ws.Cells.Item(27, 6).Value = 910
$ws.Cells.Item(28, 6).Value = 280
$ws.Cells.Item(29, 6).Value = 2195
$ws.Cells.Item(33, 6).Value = 2273
$ws.Cells.Item(34, 6).Value = 667
$ws.Cells.Item(35, 6).Value = 445
$ws.Cells.Item(36, 6).Value = 870
$ws.Cells.Item(37, 6).Value = 391
$ws.Cells.Item(38, 6).Value = 1043
$ws.Cells.Item(39, 6).Value = 864
$ws.Cells.Item(40, 6).Value = 1133
$ws.Cells.Item(41, 6).Value = 302
$ws.Cells.Item(42, 6).Value = 495
$ws.Cells.Item(44, 6).Value = 345
$ws.Cells.Item(48, 6).Value = 3441

